$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "인공지능 글 밖에 없는 블로그에 제어 글 쓰기 : 제어 엔지니어 기초"
$ws.Range("E28").Value = "https://ropiens.tistory.com/92"

$ws.Range("D37").Value = "[paper Review] GNNExplainer: Generating Explanations for Graph Neural Networks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1443&mod=document&pageid=1"

$ws.Range("D39").Value = "How to Use Normal Distribution like You Know What You Are Doing"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/How-to-Use-Normal-Distribution-like-You-Know-What-You-Are-Doing-1"

$ws.Range("D51").Value = "[북어게인 프로젝트] 헤헤부부님께서 과학고 시절을 추억하며 쓰신 ""과고라고라"""
$ws.Range("E51").Value = "https://bskyvision.com/1043"
